$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46070.01041666666
$ws.Range("B2").Value = 610.155
$ws.Range("A3").Value = 46070.02083333334
$ws.Range("B3").Value = 622.7329999999999
$ws.Range("A4").Value = 46070.03125
$ws.Range("B4").Value = 634.5650000000001
$ws.Range("A5").Value = 46070.04166666666
$ws.Range("B5").Value = 647.2910000000001
$ws.Range("A6").Value = 46070.05208333334
$ws.Range("B6").Value = 676.96
$ws.Range("A7").Value = 46070.0625
$ws.Range("B7").Value = 699.244
$ws.Range("A8").Value = 46070.07291666666
$ws.Range("B8").Value = 722.215
$ws.Range("A9").Value = 46070.08333333334
$ws.Range("B9").Value = 745.693
$ws.Range("A10").Value = 46070.09375
$ws.Range("B10").Value = 779.1799999999999
$ws.Range("A11").Value = 46070.10416666666
$ws.Range("B11").Value = 801.361
$ws.Range("A12").Value = 46070.11458333334
$ws.Range("B12").Value = 823.302
$ws.Range("A13").Value = 46070.125
$ws.Range("B13").Value = 844.577
$ws.Range("A14").Value = 46070.13541666666
$ws.Range("B14").Value = 894.139
$ws.Range("A15").Value = 46070.14583333334
$ws.Range("B15").Value = 921.274
$ws.Range("A16").Value = 46070.15625
$ws.Range("B16").Value = 949.04
$ws.Range("A17").Value = 46070.16666666666
$ws.Range("B17").Value = 976.562
$ws.Range("A18").Value = 46070.17708333334
$ws.Range("B18").Value = 1032.549
$ws.Range("A19").Value = 46070.1875
$ws.Range("B19").Value = 1072.03
$ws.Range("A20").Value = 46070.19791666666
$ws.Range("B20").Value = 1111.022
$ws.Range("A21").Value = 46070.20833333334
$ws.Range("B21").Value = 1150.906
$ws.Range("A22").Value = 46070.21875
$ws.Range("B22").Value = 1219.499
$ws.Range("A23").Value = 46070.22916666666
$ws.Range("B23").Value = 1267.366
$ws.Range("A24").Value = 46070.23958333334
$ws.Range("B24").Value = 1314.879
$ws.Range("A25").Value = 46070.25
$ws.Range("B25").Value = 1396.655
$ws.Range("A26").Value = 46070.26041666666
$ws.Range("B26").Value = 1422.915
$ws.Range("A27").Value = 46070.27083333334
$ws.Range("B27").Value = 1456.81
$ws.Range("A28").Value = 46070.28125
$ws.Range("B28").Value = 1522.308
$ws.Range("A29").Value = 46070.29166666666
$ws.Range("B29").Value = 1558.588
$ws.Range("A30").Value = 46070.30208333334
$ws.Range("B30").Value = 1589.479
$ws.Range("A31").Value = 46070.3125
$ws.Range("B31").Value = 1606.581
$ws.Range("A32").Value = 46070.32291666666
$ws.Range("B32").Value = 1620.897
$ws.Range("A33").Value = 46070.33333333334
$ws.Range("B33").Value = 1636.54
$ws.Range("A34").Value = 46070.34375
$ws.Range("B34").Value = 1643.287
$ws.Range("A35").Value = 46070.35416666666
$ws.Range("B35").Value = 1668.963
$ws.Range("A36").Value = 46070.36458333334
$ws.Range("B36").Value = 1694.65
$ws.Range("A37").Value = 46070.375
$ws.Range("B37").Value = 1721.997
$ws.Range("A38").Value = 46070.38541666666
$ws.Range("B38").Value = 1771.669
$ws.Range("A39").Value = 46070.39583333334
$ws.Range("B39").Value = 1796.292
$ws.Range("A40").Value = 46070.40625
$ws.Range("B40").Value = 1819.563
$ws.Range("A41").Value = 46070.41666666666
$ws.Range("B41").Value = 1844.172
$ws.Range("A42").Value = 46070.42708333334
$ws.Range("B42").Value = 1890.347
$ws.Range("A43").Value = 46070.4375
$ws.Range("B43").Value = 1916.943
$ws.Range("A44").Value = 46070.44791666666
$ws.Range("B44").Value = 1944.919
$ws.Range("A45").Value = 46070.45833333334
$ws.Range("B45").Value = 1974.492
$ws.Range("A46").Value = 46070.46875
$ws.Range("B46").Value = 2012.855
$ws.Range("A47").Value = 46070.47916666666
$ws.Range("B47").Value = 2036.77
$ws.Range("A48").Value = 46070.48958333334
$ws.Range("B48").Value = 2059.049
$ws.Range("A49").Value = 46070.5
$ws.Range("B49").Value = 2079.852
$ws.Range("A50").Value = 46070.51041666666
$ws.Range("B50").Value = 2110.469
$ws.Range("A51").Value = 46070.52083333334
$ws.Range("B51").Value = 2133.3
$ws.Range("A52").Value = 46070.53125
$ws.Range("B52").Value = 2141.844
$ws.Range("A53").Value = 46070.54166666666
$ws.Range("B53").Value = 2166.49
$ws.Range("A54").Value = 46070.55208333334
$ws.Range("B54").Value = 2200.902
$ws.Range("A55").Value = 46070.5625
$ws.Range("B55").Value = 2218.604
$ws.Range("A56").Value = 46070.57291666666
$ws.Range("B56").Value = 2232.56
$ws.Range("A57").Value = 46070.58333333334
$ws.Range("B57").Value = 2247.999
$ws.Range("A58").Value = 46070.59375
$ws.Range("B58").Value = 2290.145
$ws.Range("A59").Value = 46070.60416666666
$ws.Range("B59").Value = 2303.524
$ws.Range("A60").Value = 46070.61458333334
$ws.Range("B60").Value = 2316.275
$ws.Range("A61").Value = 46070.625
$ws.Range("B61").Value = 2329.355
$ws.Range("A62").Value = 46070.63541666666
$ws.Range("B62").Value = 2355.24
$ws.Range("A63").Value = 46070.64583333334
$ws.Range("B63").Value = 2363.817
$ws.Range("A64").Value = 46070.65625
$ws.Range("B64").Value = 2372.477
$ws.Range("A65").Value = 46070.66666666666
$ws.Range("B65").Value = 2314.155
$ws.Range("A66").Value = 46070.67708333334
$ws.Range("B66").Value = 2322.99
$ws.Range("A67").Value = 46070.6875
$ws.Range("B67").Value = 2398.713
$ws.Range("A68").Value = 46070.69791666666
$ws.Range("B68").Value = 2406.908
$ws.Range("A69").Value = 46070.70833333334
$ws.Range("B69").Value = 2413.852
$ws.Range("A70").Value = 46070.71875
$ws.Range("B70").Value = 2433.293
$ws.Range("A71").Value = 46070.72916666666
$ws.Range("B71").Value = 2432.504
$ws.Range("A72").Value = 46070.73958333334
$ws.Range("B72").Value = 2432.357
$ws.Range("A73").Value = 46070.75
$ws.Range("B73").Value = 2430.213
$ws.Range("A74").Value = 46070.76041666666
$ws.Range("B74").Value = 2425.504
$ws.Range("A75").Value = 46070.77083333334
$ws.Range("B75").Value = 2423.982
$ws.Range("A76").Value = 46070.78125
$ws.Range("B76").Value = 2423.517
$ws.Range("A77").Value = 46070.79166666666
$ws.Range("B77").Value = 2422.766
$ws.Range("A78").Value = 46070.80208333334
$ws.Range("B78").Value = 2426.662
$ws.Range("A79").Value = 46070.8125
$ws.Range("B79").Value = 2426.54
$ws.Range("A80").Value = 46070.82291666666
$ws.Range("B80").Value = 2426.833
$ws.Range("A81").Value = 46070.83333333334
$ws.Range("B81").Value = 2425.379
$ws.Range("A82").Value = 46070.84375
$ws.Range("B82").Value = 2418.536
$ws.Range("A83").Value = 46070.85416666666
$ws.Range("B83").Value = 2416.146
$ws.Range("A84").Value = 46070.86458333334
$ws.Range("B84").Value = 2413.646
$ws.Range("A85").Value = 46070.875
$ws.Range("B85").Value = 2411.294
$ws.Range("A86").Value = 46070.88541666666
$ws.Range("B86").Value = 2401.289
$ws.Range("A87").Value = 46070.89583333334
$ws.Range("B87").Value = 2332.134
$ws.Range("A88").Value = 46070.90625
$ws.Range("B88").Value = 2330.047
$ws.Range("A89").Value = 46070.91666666666
$ws.Range("B89").Value = 2329.357
$ws.Range("A90").Value = 46070.92708333334
$ws.Range("B90").Value = 2399.196
$ws.Range("A91").Value = 46070.9375
$ws.Range("B91").Value = 2341.932
$ws.Range("A92").Value = 46070.94791666666
$ws.Range("B92").Value = 2348.17
$ws.Range("A93").Value = 46070.95833333334
$ws.Range("B93").Value = 2354.198
$ws.Range("A94").Value = 46070.96875
$ws.Range("B94").Value = 0
$ws.Range("A95").Value = 46070.97916666666
$ws.Range("B95").Value = 0
$ws.Range("A96").Value = 46070.98958333334
$ws.Range("B96").Value = 0
$ws.Range("A97").Value = 46071
$ws.Range("B97").Value = 0
